# Penalty Reward System (unfinished) - update forecast week dates/values
# and refresh the rolled-up Summary metrics.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Forecast Comparison": shift every week forward by one week and
# refresh the MyForecast (column D) figures.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

$weekUpdates = @(
    @{ Row = 2;  Date = "2025-01-12"; Forecast = 107 },
    @{ Row = 3;  Date = "2025-01-19"; Forecast = 121 },
    @{ Row = 4;  Date = "2025-01-26"; Forecast = 148 },
    @{ Row = 5;  Date = "2025-02-02"; Forecast = 165 },
    @{ Row = 6;  Date = "2025-02-09"; Forecast = 160 },
    @{ Row = 7;  Date = "2025-02-16"; Forecast = 143 },
    @{ Row = 8;  Date = "2025-02-23"; Forecast = 140 },
    @{ Row = 9;  Date = "2025-03-02"; Forecast = 153 },
    @{ Row = 10; Date = "2025-03-09"; Forecast = 170 },
    @{ Row = 11; Date = "2025-03-16"; Forecast = 174 },
    @{ Row = 12; Date = "2025-03-23"; Forecast = 157 },
    @{ Row = 13; Date = "2025-03-30"; Forecast = 137 },
    @{ Row = 14; Date = "2025-04-06"; Forecast = 132 },
    @{ Row = 15; Date = "2025-04-13"; Forecast = 141 },
    @{ Row = 16; Date = "2025-04-20"; Forecast = 153 },
    @{ Row = 17; Date = "2025-04-27"; Forecast = 150 }
)

foreach ($u in $weekUpdates) {
    $dateCell = $ws1.Cells.Item($u.Row, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $u.Date

    $ws1.Cells.Item($u.Row, 4).Value = $u.Forecast
}

# ---------------------------------------------------------------------
# Sheet "Summary": refresh the computed roll-up metrics to match the
# shifted forecast window.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Summary")

$summaryUpdates = @(
    @{ Row = 2;  Value = "2023-01-01 to 2025-01-05" },
    @{ Row = 4;  Value = "413" },
    @{ Row = 8;  Value = "13448 units" },
    @{ Row = 9;  Value = "2351" },
    @{ Row = 10; Value = "1137" },
    @{ Row = 11; Value = "541" },
    @{ Row = 12; Value = "174" },
    @{ Row = 14; Value = "107" },
    @{ Row = 15; Value = "2025-01-12" }
)

foreach ($u in $summaryUpdates) {
    $cell = $ws2.Cells.Item($u.Row, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
}
